$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename column headers in row 1 to carry the format-version suffix instead
#    of the generic "_old" / "_new" markers:
#      <name>_old -> <name>_FV2310
#      <name>_new -> <name>_FV2404
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value2
    if ($text -ne $null) {
        $newText = $text -replace '_old$', '_FV2310'
        $newText = $newText -replace '_new$', '_FV2404'
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}

# 2) Convert the used range A1:U58 into a banded Excel Table so the new
#    headers show up as the table's column names (autofilter included).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium2"
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false

# 3) Freeze the header row: split below row 1, keep the scrolling pane's
#    top-left cell at A2.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Headers: $($ws.Cells.Item(1,1).Value2) | $($ws.Cells.Item(1,11).Value2) | $($ws.Cells.Item(1,12).Value2)"
Write-Output "Table: $($tbl.Name) ref=$($tbl.Range.Address())"
